# Exempt process emissions from carbon tax: flip the Boolean control-lever
# value on the "BEPEfCT" sheet (cell B2) from 0 (not exempt) to 1 (exempt),
# and switch the active sheet/selection from "About" to "BEPEfCT" (matching
# the author having been working on that sheet when the workbook was saved).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BEPEfCT")

# Flip the boolean lever: 0 -> 1 (exempt process emissions from carbon tax)
$ws.Range("B2").Value = 1

# Make BEPEfCT the active/selected sheet, with B3 selected, matching the
# saved workbook view state captured in the diff.
$ws.Activate()
$ws.Range("B3").Select()
